$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$sh.Table.ApplyStyle("{065B46CB-9581-4657-B4C4-BCB99A0E4069}")
